$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update retailer details and phone number for Rajshahi zone (Biswas Telecom)
$ws.Range("B2").Value = "RET-07880"
$ws.Range("C2").Value = "Biswas Telecom"
$ws.Range("D2").Value = "Nikhil Chandro Biswas"
$ws.Range("F2").Value = 1723656356

# Setting a numeric value resets the cell's "number stored as text" style;
# restore the original quote-prefixed number style by pasting formats from
# an empty cell further down that already carries it.
$ws.Range("F6").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# Rows 3-5: clear the remaining retailer rows (data removed)
$ws.Range("A3:F3").ClearContents()
$ws.Range("A4:F4").ClearContents()
$ws.Range("A5:G5").ClearContents()

# Widen column D slightly (new retailer name "Nikhil Chandro Biswas" is longer)
$ws.Range("D1").ColumnWidth = 20

# Update the active selection cell
$ws.Range("F19").Select()
